$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-autofit the data columns (B:F) to their content now that the sheet
# has been reviewed -- mirrors a user selecting the columns and choosing
# Home > Format > AutoFit Column Width (values match the widths Excel's
# AutoFit computed for this font/content).
$ws.Columns(2).ColumnWidth = 22.6
$ws.Columns(3).ColumnWidth = 13.6
$ws.Columns(4).ColumnWidth = 8.6
$ws.Columns(5).ColumnWidth = 16.6
$ws.Columns(6).ColumnWidth = 61.6

# Tighten the row heights for the data rows back down (they had been left
# very tall from manual wrapping) and mark them as explicit/custom heights.
$ws.Rows(2).RowHeight = 18
$ws.Rows(3).RowHeight = 17.25
$ws.Rows(4).RowHeight = 17.25
$ws.Rows(5).RowHeight = 17.25
$ws.Rows(6).RowHeight = 15.75

# Force the data range (including the new blank row below it) to be
# stored/displayed as Text so values like tool numbers keep any leading
# zeros / formatting verbatim.
$ws.Range("A2:G7").NumberFormat = "@"

# Move the active selection to C4, matching where the user ended up.
$ws.Range("C4").Select() | Out-Null
